# Turn the "Dummy File" roster sheet into a form template:
#  - header "Fees Paid" -> "Fees Paid ?", drop "Cutoff Cleared" header text
#  - insert a new row 2 describing each column's form-field type
#    (text / date / select / date / date / text / checkbox)
#  - the old department rows slide down to rows 3-6 and lose their
#    now-unused date values (kept as blank date-formatted cells)
#  - cosmetic: widen column C (department names) and H, and leave the
#    active selection on H2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "field type" row; existing rows 2-5 shift to 3-6
# and inherit row 1's formatting for the newly inserted row, same as
# Excel's native Insert behaviour.
$ws.Rows("2:2").Insert()

# Header row tweaks
$ws.Range("G1").Value = "Fees Paid ?"
$ws.Range("H1").Value = $null

# New row describing each field's input type
$ws.Range("A2").Value = "text"
$ws.Range("B2").Value = "date"
$ws.Range("C2").Value = "select"
$ws.Range("D2").Value = "date"
$ws.Range("E2").Value = "date"
$ws.Range("F2").Value = "text"
$ws.Range("G2").Value = "checkbox"

# The department rows (now 3-6) no longer carry sample dates
$ws.Range("D3:D6").ClearContents()

# Cosmetic column widths so the department names / checkbox column fit
$ws.Columns("C").ColumnWidth = 29.25
$ws.Columns("H").ColumnWidth = 13.59

# Leave the selection where the author left it
$ws.Range("H2").Select() | Out-Null
